$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold numeric-looking text (e.g. "42.771.37",
# "  +0.07%  "); force Text format first so Excel keeps them as strings instead
# of auto-converting to numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

$ws.Range('D2').Value = '42.771.37'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '2.571.14'
$ws.Range('E3').Value = '  +2.32%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '301.88'
$ws.Range('E5').Value = '  +3.16%  '
$ws.Range('D6').Value = '96.44'
$ws.Range('E6').Value = '  +4.38%  '
$ws.Range('D7').Value = '0.573'
$ws.Range('E7').Value = '  +1.15%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.546'
$ws.Range('E9').Value = '  +1.43%  '
$ws.Range('D10').Value = '36.13'
$ws.Range('E10').Value = '  +2.20%  '
$ws.Range('D11').Value = '0.0806'
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('D12').Value = '7.68'
$ws.Range('E12').Value = '  +1.84%  '
$ws.Range('D13').Value = '0.115'
$ws.Range('E13').Value = '  +7.88%  '
$ws.Range('D14').Value = '2.602.40'
$ws.Range('E14').Value = '  +3.21%  '
$ws.Range('D15').Value = '0.880'
$ws.Range('E15').Value = '  +3.35%  '
$ws.Range('D16').Value = '14.33'
$ws.Range('E16').Value = '  +3.07%  '
$ws.Range('D17').Value = '42.782.68'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.0₃0983'
$ws.Range('E18').Value = '  +3.51%  '
$ws.Range('B19').Value = 'InternetComputer(DFINITY)'
$ws.Range('C19').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D19').Value = '12.71'
$ws.Range('E19').Value = '  +6.23%  '
$ws.Range('D20').Value = '6.60'
$ws.Range('E20').Value = '  +2.08%  '
$ws.Range('D21').Value = '72.07'
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').Value = '253.64'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').Value = '2.94'
$ws.Range('E23').Value = '  +3.59%  '
$ws.Range('D24').Value = '2.10'
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('D25').Value = '28.17'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').Value = '10.17'
$ws.Range('E27').Value = '  +4.40%  '
$ws.Range('D28').Value = '38.85'
$ws.Range('E28').Value = '  +10.37%  '
$ws.Range('D29').Value = '2.10'
$ws.Range('E29').Value = '  -4.99%  '
$ws.Range('D30').Value = '5.97'
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('D31').Value = '155.36'
$ws.Range('E31').Value = '  +4.04%  '
$ws.Range('D32').Value = '2.16'
$ws.Range('E32').Value = '  +2.01%  '
$ws.Range('D33').Value = '2.75'
$ws.Range('E33').Value = '  +1.58%  '
$ws.Range('D34').Value = '0.0805'
$ws.Range('E34').Value = '  +2.74%  '
$ws.Range('D35').Value = '3.33'
$ws.Range('E35').Value = '  -1.82%  '
$ws.Range('D36').Value = '18.19'
$ws.Range('E36').Value = '  +15.45%  '
$ws.Range('D37').Value = '0.113'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '0.119'
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('D39').Value = '23.51'
$ws.Range('E39').Value = '  +3.19%  '
$ws.Range('D40').Value = '2.13'
$ws.Range('E40').Value = '  +34.08%  '
$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D41').Value = '3.39'
$ws.Range('E41').Value = '  +1.52%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '3.84'
$ws.Range('E42').Value = '  +3.26%  '
$ws.Range('D43').Value = '0.0306'
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('D44').Value = '2.079.48'
$ws.Range('E44').Value = '  +2.54%  '
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').Value = '9.23'
$ws.Range('E46').Value = '  +5.19%  '
$ws.Range('D47').Value = '84.56'
$ws.Range('E47').Value = '  +1.42%  '
$ws.Range('D48').Value = '77.15'
$ws.Range('E48').Value = '  +14.22%  '
$ws.Range('D49').Value = '2.817.76'
$ws.Range('E49').Value = '  +2.17%  '
$ws.Range('D50').Value = '105.09'
$ws.Range('E50').Value = '  +3.47%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.189'
$ws.Range('E51').Value = '  +3.26%  '
